$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fix error in Device mapping: the zib mapping data (columns B-H,
# "MedicalDevice.Product.ProductType") was incorrectly entered on row 8
# (EHDSDevice.modelNumber) instead of row 12 (EHDSDevice.type).
# Move that data from row 8 to row 12.

$values = @(
    "MedicalDevice.Product.ProductType",
    "ProductType",
    "MedicalDevice.Product.ProductType",
    "CD",
    "0..1",
    "NL-CM:10.1.3",
    "The code of the type of product."
)

# Clear row 8, columns B:H (the misplaced data)
$ws.Range("B8:H8").Value = ""

# Write the data into row 12, columns B:H (the correct location)
for ($col = 2; $col -le 8; $col++) {
    $ws.Cells.Item(12, $col).Value = $values[$col - 2]
}
